$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Cthrc1"
$ws.Cells.Item(2, 3).Value = "Ror2"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.01838633333333333
$ws.Cells.Item(2, 8).Value = 0.055159
$ws.Cells.Item(2, 9).Value = 0.003339500866342531
$ws.Cells.Item(2, 10).Value = 0.003339500866342531
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.294804333333333
$ws.Cells.Item(2, 14).Value = 6.884412999999999
$ws.Cells.Item(2, 15).Value = 0.6595351916216082
$ws.Cells.Item(2, 16).Value = 0.6595351916216082
$ws.Cells.Item(2, 17).Value = 0.04219303740744444
$ws.Cells.Item(2, 18).Value = 0.379737336667
$ws.Cells.Item(2, 19).Value = 0.002202518343803748
$ws.Cells.Item(2, 20).Value = 0.002202518343803748

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Cthrc1"
$ws.Cells.Item(3, 3).Value = "Ror2"
$ws.Cells.Item(3, 4).Value = "MuSCs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.01838633333333333
$ws.Cells.Item(3, 8).Value = 0.055159
$ws.Cells.Item(3, 9).Value = 0.003339500866342531
$ws.Cells.Item(3, 10).Value = 0.003339500866342531
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.184622333333333
$ws.Cells.Item(3, 14).Value = 3.553867
$ws.Cells.Item(3, 15).Value = 0.3404648083783919
$ws.Cells.Item(3, 16).Value = 0.3404648083783919
$ws.Cells.Item(3, 17).Value = 0.02178086109477778
$ws.Cells.Item(3, 18).Value = 0.196027749853
$ws.Cells.Item(3, 19).Value = 0.001136982522538784
$ws.Cells.Item(3, 20).Value = 0.001136982522538783

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Cthrc1"
$ws.Cells.Item(4, 3).Value = "Ror2"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 5.449420333333333
$ws.Cells.Item(4, 8).Value = 16.348261
$ws.Cells.Item(4, 9).Value = 0.9897755900704113
$ws.Cells.Item(4, 10).Value = 0.9897755900704112
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.294804333333333
$ws.Cells.Item(4, 14).Value = 6.884412999999999
$ws.Cells.Item(4, 15).Value = 0.6595351916216082
$ws.Cells.Item(4, 16).Value = 0.6595351916216082
$ws.Cells.Item(4, 17).Value = 12.50535339508811
$ws.Cells.Item(4, 18).Value = 112.548180555793
$ws.Cells.Item(4, 19).Value = 0.652791833459479
$ws.Cells.Item(4, 20).Value = 0.652791833459479

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Cthrc1"
$ws.Cells.Item(5, 3).Value = "Ror2"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 5.449420333333333
$ws.Cells.Item(5, 8).Value = 16.348261
$ws.Cells.Item(5, 9).Value = 0.9897755900704113
$ws.Cells.Item(5, 10).Value = 0.9897755900704112
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.184622333333333
$ws.Cells.Item(5, 14).Value = 3.553867
$ws.Cells.Item(5, 15).Value = 0.3404648083783919
$ws.Cells.Item(5, 16).Value = 0.3404648083783919
$ws.Cells.Item(5, 17).Value = 6.455505030587444
$ws.Cells.Item(5, 18).Value = 58.099545275287
$ws.Cells.Item(5, 19).Value = 0.3369837566109323
$ws.Cells.Item(5, 20).Value = 0.3369837566109323

# Row 6
$ws.Cells.Item(6, 1).Value = "MuSCs"
$ws.Cells.Item(6, 2).Value = "Cthrc1"
$ws.Cells.Item(6, 3).Value = "Ror2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.03790633333333333
$ws.Cells.Item(6, 8).Value = 0.113719
$ws.Cells.Item(6, 9).Value = 0.006884909063246366
$ws.Cells.Item(6, 10).Value = 0.006884909063246365
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.294804333333333
$ws.Cells.Item(6, 14).Value = 6.884412999999999
$ws.Cells.Item(6, 15).Value = 0.6595351916216082
$ws.Cells.Item(6, 16).Value = 0.6595351916216082
$ws.Cells.Item(6, 17).Value = 0.0869876179941111
$ws.Cells.Item(6, 18).Value = 0.782888561947
$ws.Cells.Item(6, 19).Value = 0.00454083981832554
$ws.Cells.Item(6, 20).Value = 0.004540839818325538

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Cthrc1"
$ws.Cells.Item(7, 3).Value = "Ror2"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.03790633333333333
$ws.Cells.Item(7, 8).Value = 0.113719
$ws.Cells.Item(7, 9).Value = 0.006884909063246366
$ws.Cells.Item(7, 10).Value = 0.006884909063246365
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.184622333333333
$ws.Cells.Item(7, 14).Value = 3.553867
$ws.Cells.Item(7, 15).Value = 0.3404648083783919
$ws.Cells.Item(7, 16).Value = 0.3404648083783919
$ws.Cells.Item(7, 17).Value = 0.04490468904144444
$ws.Cells.Item(7, 18).Value = 0.404142201373
$ws.Cells.Item(7, 19).Value = 0.002344069244920828
$ws.Cells.Item(7, 20).Value = 0.002344069244920827

